$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "linear_product_length_out_of_store"
$ws.Range("B3").Value = "sub_category"
$ws.Range("B3").Font.Name = "Calibri"
$ws.Range("B3").Font.Color = 0

$ws.Range("C3").Value = "COLD & FLU ; IBS ; KIDS COLD & FLU ; KIDS COUGH ;  KIDS DECONGESTANTS ; KIDS DIGESTIVE HEALTH ;  KIDS HAYFEVER ; LAXATIVES ; KIDS TEETHING ; WIND ; ALLERGY ; REHYDRATION ; PROBITOTICS ; PAIN MANAGEMENT ; FIRST AID ; MIGRAINE RELIEF ; RASH TREATMENT ; ANTI-AGE FACE ; ARTIFICIAL TAN ; BODY CLEANSING ; COSMETICS ; DEODORANTS ; MEN'S TOILETRIES ; SUNCARE ; HAIR CARE ; BABY HEALTHCARE ; FOR MUM ; KIDS HAIRCARE ; KIDS TOILETRIES ; KIDS WIPES ; BABY SUNCARE ; COTTON ; INCONTINENCE ; SANITARY TOWELS ; FEMININE WASH ; KIDS MOUTHWASH"
$ws.Range("C3").Font.Name = "Arial"
$ws.Range("C3").Font.Color = 3355443

$ws.Range("C9").Select()
